$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (subject numbers)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) values - C2 and E2 are cleared entirely
$ws.Range("B2").Value = 7.9032871048677933
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 8.1534275168136787
$ws.Range("E2").ClearContents()

# Row 3 (STR) values
$ws.Range("B3").Value = 5.9304623509933521
$ws.Range("C3").Value = -6.240402481620599
$ws.Range("D3").Value = 6.1208698948263134
$ws.Range("E3").Value = -8.9813355715450243

# Update the selection to match the new active range
$ws.Range("B1:E3").Select()
